$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.941.20"
$ws.Range("E2").Value = "  -3.34%  "
$ws.Range("D3").Value = "1.719.90"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "309.71"
$ws.Range("E5").Value = "  -5.88%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.4872"
$ws.Range("E7").Value = "  +7.37%  "
$ws.Range("D8").Value = "0.3492"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "42.09"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "0.07261"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "1.046"
$ws.Range("E11").Value = "  -4.53%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "19.90"
$ws.Range("E13").Value = "  -3.91%  "
$ws.Range("D14").Value = "5.862"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "1.727.32"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").Value = "6.851"
$ws.Range("E16").Value = "  -4.69%  "
$ws.Range("D17").Value = "86.69"
$ws.Range("E17").Value = "  -6.00%  "
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").Value = "0.06375"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D21").Value = "16.52"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").Value = "5.644"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").Value = "26.998.70"
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("D24").Value = "10.81"
$ws.Range("E24").Value = "  -3.79%  "
$ws.Range("D25").Value = "2.083"
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("D26").Value = "154.01"
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("D27").Value = "19.92"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "1.918.12"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").Value = "2.066"
$ws.Range("E29").Value = "  -4.74%  "
$ws.Range("D30").Value = "120.75"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("E31").Value = "  -4.64%  "
$ws.Range("D32").Value = "0.09299"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").Value = "5.357"
$ws.Range("E34").Value = "  -4.01%  "
$ws.Range("D35").Value = "0.05900"
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("D36").Value = "0.02177"
$ws.Range("E36").Value = "  -4.18%  "
$ws.Range("D37").Value = "1.440"
$ws.Range("E37").Value = "  +4.37%  "
$ws.Range("D38").Value = "10.98"
$ws.Range("E38").Value = "  -6.94%  "
$ws.Range("D39").Value = "0.1996"
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("D40").Value = "4.745"
$ws.Range("E40").Value = "  -3.89%  "
$ws.Range("D42").Value = "0.5989"
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("D43").Value = "1.088"
$ws.Range("E43").Value = "  -7.97%  "
$ws.Range("D44").Value = "7.510"
$ws.Range("E44").Value = "  -4.16%  "
$ws.Range("D45").Value = "12.75"
$ws.Range("E45").Value = "  -3.32%  "
$ws.Range("D46").Value = "3.576"
$ws.Range("E46").Value = "  -4.23%  "
$ws.Range("D47").Value = "0.5623"
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "117.83"
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("D49").Value = "1.835"
$ws.Range("E49").Value = "  -5.22%  "
$ws.Range("D50").Value = "1.111"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("E51").Value = "  -2.41%  "
